# Apply odds updates to the FlashScore weekly games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value  = 1.75
$ws.Range("H2").Value  = 3.5
$ws.Range("M2").Value  = 1.07
$ws.Range("N2").Value  = 9
$ws.Range("R2").Value  = 1.58
$ws.Range("AH2").Value = 11
$ws.Range("AO2").Value = 9.5
$ws.Range("AR2").Value = 51
$ws.Range("AU2").Value = 9

# Row 3 updates
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
